# Completed HashTable, Store, Hash and DataFile classes
#
# The worksheet used to have a generic header row
# (category, product, expiratin_date, weight, cost) styled with a
# centered/shaded format. It is replaced with a header row that lists
# the distinct category names found in column A (fruits, vegetables,
# ice-creams, water, soda, sweets) in plain/default formatting, which
# also extends the used range by one column (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the special centered/shaded formatting the old header row used
# so the header cells fall back to the default style.
$ws.Range("A1:E1").Style = "Normal"

# Write the new header labels (category names), adding a 6th column.
$ws.Range("A1").Value = "fruits"
$ws.Range("B1").Value = "vegetables"
$ws.Range("C1").Value = "ice-creams"
$ws.Range("D1").Value = "water"
$ws.Range("E1").Value = "soda"
$ws.Range("F1").Value = "sweets"

# Match the author's final selection.
$ws.Range("I14").Select()
